# feat: add 2022-Q4 data
#
# The workbook currently has 3 sheets: "总计" (totals), "2022-Q3", "2022-Q1".
# We add a new "2022-Q4" sheet (a copy of the "2022-Q3" sheet's layout,
# populated with the new quarter's fund-holding figures) positioned right
# after "总计" and before "2022-Q3", and we add a corresponding summary row
# to the "总计" sheet.

$wb = $excel.ActiveWorkbook

$wsTotal = $wb.Worksheets.Item("总计")
$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# --- 1. Create the new "2022-Q4" sheet by duplicating "2022-Q3" and
#        placing the duplicate immediately before it. -------------------
$wsQ3.Copy($wsQ3)
$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# Update the duplicated rows with the 2022-Q4 figures (note the fund
# ordering is swapped relative to the 2022-Q3 sheet: the USD-denominated
# share class now comes first).
$wsQ4.Range("B2").Value = "'000927"
$wsQ4.Range("C2").Value = "博时大中华亚太精选股票（QDII）美元现汇"
$wsQ4.Range("D2").Value = "'0.28"
$wsQ4.Range("E2").Value = "'88.48"
$wsQ4.Range("F2").Value = "'4.02"
$wsQ4.Range("G2").Value = "'0.0113"
$wsQ4.Range("H2").Value = 10

$wsQ4.Range("B3").Value = "'050015"
$wsQ4.Range("C3").Value = "博时大中华亚太精选股票（QDII）人民币"
$wsQ4.Range("D3").Value = "'0.28"
$wsQ4.Range("E3").Value = "'88.48"
$wsQ4.Range("F3").Value = "'4.02"
$wsQ4.Range("G3").Value = "'0.0113"
$wsQ4.Range("H3").Value = 10

# --- 2. Add a 2022-Q4 row to the "总计" summary sheet, pushing the
#        existing 2022-Q3 / 2022-Q1 rows down by one. --------------------

# Extend the index-column styling (column A) down to the new row 4.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A4").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "'2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.02

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "'2022-Q3"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.02

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "'2022-Q1"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 0.05
